$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    try {
      $val = $cell.Value2
      $new = $val.Replace("D64","D69").Replace("D51","D55").Replace("D80","D86").Replace("S30","S31")
      if ($new -ne $val) {
        $cell.Value2 = $new
      }
    } catch {
    }
  }
}
